$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 321:322, shifting existing rows 321-424 down to 323-426
$ws.Rows("321:322").Insert()

# Row 321: new Copenhague / Primera record
$ws.Range("A321").Value = 7
$ws.Range("B321").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C321").Value = "Ñuble"
$ws.Range("D321").Value = 45093
$ws.Range("E321").Value = 16
$ws.Range("F321").Value = 100112006
$ws.Range("G321").Value = "Repollo"
$ws.Range("H321").Value = "Copenhague"
$ws.Range("I321").Value = "Primera"
$ws.Range("J321").Value = 400
$ws.Range("K321").Value = 1200
$ws.Range("L321").Value = 1200
$ws.Range("M321").Value = 1200
$ws.Range("N321").Value = "$/unidad"
$ws.Range("O321").Value = "Provincia de Diguillín"
$ws.Range("P321").Value = 1200
$ws.Range("Q321").Value = 1
$ws.Range("R321").Value = "Hortaliza"

# Row 322: new Copenhague / Segunda record
$ws.Range("A322").Value = 7
$ws.Range("B322").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C322").Value = "Ñuble"
$ws.Range("D322").Value = 45093
$ws.Range("E322").Value = 16
$ws.Range("F322").Value = 100112006
$ws.Range("G322").Value = "Repollo"
$ws.Range("H322").Value = "Copenhague"
$ws.Range("I322").Value = "Segunda"
$ws.Range("J322").Value = 300
$ws.Range("K322").Value = 1000
$ws.Range("L322").Value = 1000
$ws.Range("M322").Value = 1000
$ws.Range("N322").Value = "$/unidad"
$ws.Range("O322").Value = "Provincia de Diguillín"
$ws.Range("P322").Value = 1000
$ws.Range("Q322").Value = 1
$ws.Range("R322").Value = "Hortaliza"

"done"
